# fall 23 week 7 updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(3,0,3,3),
    @(5,1,3,2),
    @(5,2,5,1),
    @(3,2,7,0),
    @(4,0,5,3),
    @(3,2,3,1),
    @(3,1,5,2),
    @(6,0,7,2),
    @(6,2,5,1),
    @(3,3,2,0),
    @(4,0,4,2),
    @(4,2,5,0),
    @(3,2,5,0),
    @(5,2,4,0),
    @(5,2,5,1),
    @(3,0,3,3),
    @(5,2,6,0),
    @(4,2,5,1),
    @(4,2,5,0),
    @(2,2,3,0),
    @(3,2,4,0),
    @(4,1,4,2),
    @(6,0,6,3),
    @(4,1,6,2),
    @(5,2,3,1)
)

$startRow = 2318
$endRow = $startRow + $data.Length - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$ws.Application.ActiveWindow.ScrollRow = 2317
$ws.Range("A2343").Select()
